# Applies the diff to statement_74.xlsx:
#  - Cardholder name / number changes (row 2-3)
#  - Statement period changes (KONTOSTAND dates)
#  - Transaction rows 6-9 updated in place
#  - Transaction rows 10-11 cleared (fewer transactions this period)
#  - Closing balance / next billing date updated (rows 12-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell (well outside the printed A1:G14 area) used as a relay so a
# purely-numeric string (the card number) round-trips through the clipboard
# as TEXT instead of being auto-coerced to a Number by a direct
# Range.Value assignment. TEXT() always yields a string result, and
# PasteSpecial (values only) preserves that string-ness on the target cell
# while leaving its existing style/number-format untouched.
$scratch = $ws.Range("Z100")

function Set-TextValue($range, [string]$text) {
    $scratch.Formula = '=TEXT("' + $text + '","@")'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

# --- Row 2: cardholder first name ---
$ws.Range("C2").Value = "Hartmut"

# --- Row 3: card number + surname ---
Set-TextValue $ws.Range("B3") "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Row 5: opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 15.08.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "16.08."
$ws.Range("C6").Value = "17.08."
$ws.Range("D6").Value = "PAYPAL OLYBCN"
$ws.Range("E6").Value = "19,54-"

# --- Row 7 ---
$ws.Range("B7").Value = "18.08."
$ws.Range("C7").Value = "19.08."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-81549747"
$ws.Range("E7").Value = "53,43-"

# --- Row 8 ---
$ws.Range("B8").Value = "20.08."
$ws.Range("C8").Value = "21.08."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,39-"

# --- Row 9 ---
$ws.Range("B9").Value = "22.08."
$ws.Range("C9").Value = "23.08."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 10752056"
$ws.Range("E9").Value = "41,99-"

# --- Row 10: no longer used this period, clear it out ---
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# --- Row 11: no longer used this period, clear it out ---
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# --- Row 12: closing balance ---
$ws.Range("D12").Value = "KONTOSTAND AM 27.08.2024"
$ws.Range("E12").Value = "140,35-"

# --- Row 13: next billing date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 03.09.2024"

# Tidy up the scratch cell so it doesn't linger in the saved sheet.
$scratch.Value = ""
